$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2-64 need to be rotated: the last data row (64) moves to the
# top of the data block (row 2), and every other data row (2-63) shifts down
# by one row (becoming rows 3-64).
#
# Achieve this with a genuine Excel row move: copy the last row, insert it
# (with a shift) above the first data row, then remove the now-duplicated
# original row that got pushed one row further down.

$ws.Rows("64:64").Copy()
$ws.Rows("2:2").Insert()
$ws.Application.CutCopyMode = $false
$ws.Rows("65:65").Delete()
